# Apply edits described by the commit:
#  1. Update several meter-read counts on the existing "2020-08-07" sheet
#     (columns E and J, rows 30-38).
#  2. Add a new worksheet "2020-08-14" (a copy of "2020-08-07") at the end
#     of the workbook, with its own set of values in columns E and J
#     (rows 30-39).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update existing sheet "2020-08-07"
# ---------------------------------------------------------------------
$ws0807 = $wb.Worksheets.Item("2020-08-07")

$ws0807.Range("E30").Value = 1.0
$ws0807.Range("J30").Value = 1.0

$ws0807.Range("E31").Value = 0.0
$ws0807.Range("J31").Value = 0.0

$ws0807.Range("J32").Value = 1.0

$ws0807.Range("J33").Value = 1.0

$ws0807.Range("J34").Value = 0.0

$ws0807.Range("J35").Value = 1.0

$ws0807.Range("J36").Value = 7.0

$ws0807.Range("J37").Value = 50.0

$ws0807.Range("J38").Value = 29.0

# ---------------------------------------------------------------------
# 2) Add new sheet "2020-08-14" as a copy of "2020-08-07", placed after
#    the last existing sheet.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws0807.Copy($null, $lastSheet)

$ws0814 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws0814.Name = "2020-08-14"

$ws0814.Range("E30").Value = 1.0
$ws0814.Range("J30").Value = 1.0

$ws0814.Range("E31").Value = 0.0
$ws0814.Range("J31").Value = 0.0

$ws0814.Range("E32").Value = 0.0
$ws0814.Range("J32").Value = 1.0

$ws0814.Range("E33").Value = 0.0
$ws0814.Range("J33").Value = 1.0

$ws0814.Range("E34").Value = 1.0
$ws0814.Range("J34").Value = 0.0

$ws0814.Range("E35").Value = 0.0
$ws0814.Range("J35").Value = 1.0

$ws0814.Range("E36").Value = 0.0
$ws0814.Range("J36").Value = 7.0

$ws0814.Range("E37").Value = 1.0
$ws0814.Range("J37").Value = 50.0

$ws0814.Range("E38").Value = 0.0
$ws0814.Range("J38").Value = 29.0

$ws0814.Range("E39").Value = 0
$ws0814.Range("J39").Value = 11

# Restore the originally active sheet so we don't leave unrelated
# side effects on the workbook's active tab.
$wb.Worksheets.Item("2020-08-03").Activate()
